# "Quang Linh gộp code lại rồi nè" - merge a duplicate MaPM=15 record into
# the list, right before the MaPM=16 block, pushing that block down a row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record at row 36 (existing rows 36-38 shift to 37-39).
$ws.Rows.Item(36).Insert()

# Row 35 already holds MaPM "15" (shared string), so copy it down for A36
# to reuse the exact same shared-string cell type/value.
$ws.Range("A35").Copy()
$ws.Range("A36").PasteSpecial()

# Row 11 already holds MaSach "3" / SoLuong "3" (shared string) for B/C,
# copy that pair into B36:C36.
$ws.Range("B11:C11").Copy()
$ws.Range("B36:C36").PasteSpecial()

$excel.CutCopyMode = $false
